$wb = $excel.ActiveWorkbook

# Sheet: Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 11).Value = 5303
$ws.Cells.Item(3, 11).Value = 5474
$ws.Cells.Item(4, 6).Value = 1913
$ws.Cells.Item(4, 8).Value = 1742
$ws.Cells.Item(4, 11).Value = 1134
$ws.Cells.Item(5, 11).Value = 391
$ws.Cells.Item(6, 11).Value = 6077
$ws.Cells.Item(7, 6).Value = 24106
$ws.Cells.Item(7, 8).Value = 26055
$ws.Cells.Item(7, 11).Value = 18379

# Sheet: Grant Park
$ws = $wb.Worksheets.Item('Grant Park')
$ws.Cells.Item(5, 11).Value = 9
$ws.Cells.Item(6, 11).Value = 18

# Sheet: Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(2, 11).Value = 60
$ws.Cells.Item(7, 11).Value = 236

# Sheet: Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(3, 11).Value = 370
$ws.Cells.Item(6, 11).Value = 419
$ws.Cells.Item(7, 11).Value = 1234

# Sheet: South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 11).Value = 137
$ws.Cells.Item(7, 11).Value = 407

# Sheet: Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(6, 11).Value = 230
$ws.Cells.Item(7, 11).Value = 787

# Sheet: Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 11).Value = 178
$ws.Cells.Item(3, 11).Value = 205
$ws.Cells.Item(4, 11).Value = 30
$ws.Cells.Item(5, 11).Value = 28
$ws.Cells.Item(6, 11).Value = 181
$ws.Cells.Item(7, 11).Value = 622

# Sheet: New City
$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(3, 11).Value = 105
$ws.Cells.Item(7, 11).Value = 417

# Sheet: Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(3, 11).Value = 130
$ws.Cells.Item(6, 11).Value = 80
$ws.Cells.Item(7, 11).Value = 312

# Sheet: By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(7, 11).Value = 545
$ws.Cells.Item(8, 11).Value = 1234
$ws.Cells.Item(10, 11).Value = 103
$ws.Cells.Item(18, 11).Value = 124
$ws.Cells.Item(19, 11).Value = 541
$ws.Cells.Item(20, 11).Value = 425
$ws.Cells.Item(29, 11).Value = 982
$ws.Cells.Item(31, 11).Value = 201
$ws.Cells.Item(33, 11).Value = 787
$ws.Cells.Item(34, 11).Value = 105
$ws.Cells.Item(36, 11).Value = 243
$ws.Cells.Item(37, 11).Value = 622
$ws.Cells.Item(38, 11).Value = 18
$ws.Cells.Item(41, 11).Value = 128
$ws.Cells.Item(46, 11).Value = 37
$ws.Cells.Item(48, 11).Value = 232
$ws.Cells.Item(49, 11).Value = 102
$ws.Cells.Item(51, 11).Value = 230
$ws.Cells.Item(52, 11).Value = 479
$ws.Cells.Item(53, 11).Value = 236
$ws.Cells.Item(57, 11).Value = 68
$ws.Cells.Item(58, 11).Value = 11
$ws.Cells.Item(60, 11).Value = 115
$ws.Cells.Item(63, 6).Value = 197
$ws.Cells.Item(63, 8).Value = 292
$ws.Cells.Item(63, 11).Value = 50
$ws.Cells.Item(65, 11).Value = 417
$ws.Cells.Item(67, 11).Value = 695
$ws.Cells.Item(71, 11).Value = 58
$ws.Cells.Item(73, 11).Value = 158
$ws.Cells.Item(76, 11).Value = 254
$ws.Cells.Item(77, 11).Value = 129
$ws.Cells.Item(78, 11).Value = 211
$ws.Cells.Item(79, 11).Value = 457
$ws.Cells.Item(83, 11).Value = 407
$ws.Cells.Item(84, 11).Value = 136
$ws.Cells.Item(85, 11).Value = 867
$ws.Cells.Item(89, 11).Value = 268
$ws.Cells.Item(93, 11).Value = 69
$ws.Cells.Item(94, 11).Value = 243
$ws.Cells.Item(96, 11).Value = 199
$ws.Cells.Item(97, 11).Value = 147
$ws.Cells.Item(99, 11).Value = 312
$ws.Cells.Item(101, 6).Value = 24106
$ws.Cells.Item(101, 8).Value = 26055
$ws.Cells.Item(101, 11).Value = 18379

# Sheet: Gage Park
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(3, 11).Value = 49
$ws.Cells.Item(6, 11).Value = 72
$ws.Cells.Item(7, 11).Value = 201

# Sheet: North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 11).Value = 248
$ws.Cells.Item(7, 11).Value = 695

# Sheet: South Deering
$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(3, 11).Value = 54
$ws.Cells.Item(7, 11).Value = 136

# Sheet: Lincoln Park
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(2, 11).Value = 20
$ws.Cells.Item(7, 11).Value = 102

# Sheet: Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 11).Value = 281
$ws.Cells.Item(3, 11).Value = 355
$ws.Cells.Item(7, 11).Value = 982

# Sheet: Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(2, 11).Value = 32
$ws.Cells.Item(7, 11).Value = 232

# Sheet: Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(3, 11).Value = 171
$ws.Cells.Item(6, 11).Value = 172
$ws.Cells.Item(7, 11).Value = 541

# Sheet: River North
$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(4, 11).Value = 16
$ws.Cells.Item(7, 11).Value = 254

# Sheet: Hermosa
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Cells.Item(2, 11).Value = 45
$ws.Cells.Item(7, 11).Value = 128

# Sheet: Avondale
$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(2, 11).Value = 29
$ws.Cells.Item(3, 11).Value = 18
$ws.Cells.Item(7, 11).Value = 103

# Sheet: Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(6, 11).Value = 76
$ws.Cells.Item(7, 11).Value = 211

# Sheet: Jefferson Park
$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Cells.Item(3, 11).Value = 10
$ws.Cells.Item(7, 11).Value = 37

# Sheet: West Ridge
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(3, 11).Value = 38
$ws.Cells.Item(7, 11).Value = 199

# Sheet: Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 11).Value = 154
$ws.Cells.Item(3, 11).Value = 146
$ws.Cells.Item(4, 11).Value = 30
$ws.Cells.Item(6, 11).Value = 111
$ws.Cells.Item(7, 11).Value = 457

# Sheet: Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 11).Value = 140
$ws.Cells.Item(3, 11).Value = 136
$ws.Cells.Item(7, 11).Value = 425

# Sheet: Calumet Heights
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(3, 11).Value = 39
$ws.Cells.Item(7, 11).Value = 124

# Sheet: Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(3, 11).Value = 71
$ws.Cells.Item(7, 11).Value = 243

# Sheet: West Lawn
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(4, 11).Value = 5
$ws.Cells.Item(7, 11).Value = 69

# Sheet: Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 11).Value = 185
$ws.Cells.Item(3, 11).Value = 175
$ws.Cells.Item(5, 11).Value = 22
$ws.Cells.Item(6, 11).Value = 145
$ws.Cells.Item(7, 11).Value = 545

# Sheet: Garfield Ridge
$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(2, 11).Value = 37
$ws.Cells.Item(7, 11).Value = 105

# Sheet: West Loop
$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(2, 11).Value = 69
$ws.Cells.Item(6, 11).Value = 105
$ws.Cells.Item(7, 11).Value = 243

# Sheet: Portage Park
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(4, 11).Value = 11
$ws.Cells.Item(7, 11).Value = 158

# Sheet: West Town
$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(3, 11).Value = 28
$ws.Cells.Item(7, 11).Value = 147

# Sheet: Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(2, 11).Value = 75
$ws.Cells.Item(3, 11).Value = 82
$ws.Cells.Item(6, 11).Value = 80
$ws.Cells.Item(7, 11).Value = 268

# Sheet: Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(6, 11).Value = 79
$ws.Cells.Item(7, 11).Value = 230

# Sheet: Mckinley Park
$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Cells.Item(3, 11).Value = 14
$ws.Cells.Item(7, 11).Value = 68

# Sheet: Morgan Park
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(3, 11).Value = 34
$ws.Cells.Item(7, 11).Value = 115

# Sheet: South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 11).Value = 290
$ws.Cells.Item(3, 11).Value = 293
$ws.Cells.Item(6, 11).Value = 208
$ws.Cells.Item(7, 11).Value = 867

# Sheet: Oakland
$ws = $wb.Worksheets.Item('Oakland')
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(7, 11).Value = 58

# Sheet: Riverdale
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Cells.Item(2, 11).Value = 57
$ws.Cells.Item(7, 11).Value = 129

# Sheet: Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(3, 11).Value = 135
$ws.Cells.Item(6, 11).Value = 174
$ws.Cells.Item(7, 11).Value = 479

# Sheet: Millenium Park
$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Cells.Item(6, 11).Value = 9
$ws.Cells.Item(7, 11).Value = 11
